$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.999.71"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.772.31"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E3").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "629.04"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E5").NumberFormat = "General"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E6").NumberFormat = "General"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.770.23"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("E7").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E8").NumberFormat = "General"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E9").NumberFormat = "General"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.456"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("E11").NumberFormat = "General"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("E12").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("E13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.77"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("E14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.406.14"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("E15").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.765.04"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("E16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.967.92"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.65"
$ws.Range("D18").NumberFormat = "General"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E19").NumberFormat = "General"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("E20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "469.13"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("E21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.51"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("E22").NumberFormat = "General"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.05"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("E24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000140"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.13%  "
$ws.Range("E25").NumberFormat = "General"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.13"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("E26").NumberFormat = "General"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("E27").NumberFormat = "General"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("E28").NumberFormat = "General"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E29").NumberFormat = "General"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.921.37"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E30").NumberFormat = "General"

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("B31").NumberFormat = "General"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C31").NumberFormat = "General"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.26"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E31").NumberFormat = "General"

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("B32").NumberFormat = "General"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C32").NumberFormat = "General"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.66"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E32").NumberFormat = "General"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("E33").NumberFormat = "General"

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Kaspa"
$ws.Range("B34").NumberFormat = "General"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C34").NumberFormat = "General"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.176"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +18.91%  "
$ws.Range("E34").NumberFormat = "General"

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("B35").NumberFormat = "General"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C35").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.41"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E35").NumberFormat = "General"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E36").NumberFormat = "General"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.723.85"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E37").NumberFormat = "General"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.88"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("E38").NumberFormat = "General"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E39").NumberFormat = "General"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.25"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("E40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("E41").NumberFormat = "General"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("B42").NumberFormat = "General"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C42").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E42").NumberFormat = "General"

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Mantle"
$ws.Range("B43").NumberFormat = "General"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C43").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.961"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("E43").NumberFormat = "General"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.25%  "
$ws.Range("E45").NumberFormat = "General"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.53"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E46").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.40"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.96"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.41"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("E49").NumberFormat = "General"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.33"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.28%  "
$ws.Range("E51").NumberFormat = "General"

Write-Output "Applied cryptos update"
